# Auto-generated edit script applying the Malboro_Profits.xlsx numeric updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 8).Value = 115.666664   # H2 (ALC): 90.375 -> 115.666664
$ws.Cells.Item(2, 9).Value = 118.8   # I2 (ALC): 90.375 -> 118.8
$ws.Cells.Item(2, 10).Value = 100   # J2 (ALC): 0 -> 100
$ws.Cells.Item(2, 11).Value = 118.8   # K2 (ALC): 90.375 -> 118.8
$ws.Cells.Item(2, 12).Value = 100   # L2 (ALC): 0 -> 100
$ws.Cells.Item(2, 13).Value = -5.799999999999997   # M2 (ALC): 22.625 -> -5.799999999999997
$ws.Cells.Item(2, 14).Value = -326   # N2 (ALC): None -> -326
$ws.Cells.Item(5, 8).Value = 69.27273   # H5 (ALC): 94.666664 -> 69.27273
$ws.Cells.Item(5, 9).Value = 25.714285   # I5 (ALC): 109.4 -> 25.714285
$ws.Cells.Item(5, 10).Value = 145.5   # J5 (ALC): 21 -> 145.5
$ws.Cells.Item(5, 11).Value = 25.714285   # K5 (ALC): 109.4 -> 25.714285
$ws.Cells.Item(5, 12).Value = 145.5   # L5 (ALC): 21 -> 145.5
$ws.Cells.Item(5, 13).Value = 89.285715   # M5 (ALC): 5.599999999999994 -> 89.285715
$ws.Cells.Item(5, 14).Value = -375.5   # N5 (ALC): -251 -> -375.5
$ws.Cells.Item(17, 8).Value = 1366609.2   # H17 (ALC): 1466591.8 -> 1366609.2
$ws.Cells.Item(17, 10).Value = 1366609.2   # J17 (ALC): 1466591.8 -> 1366609.2
$ws.Cells.Item(17, 12).Value = 4099827.6   # L17 (ALC): 4399775.4 -> 4099827.6
$ws.Cells.Item(17, 14).Value = -4100163.6   # N17 (ALC): -4400111.4 -> -4100163.6
$ws.Cells.Item(33, 8).Value = 10419811   # H33 (ALC): 9618646 -> 10419811
$ws.Cells.Item(33, 9).Value = 17858660   # I33 (ALC): 16668361 -> 17858660
$ws.Cells.Item(33, 10).Value = 5420.2   # J33 (ALC): 5397.1816 -> 5420.2
$ws.Cells.Item(33, 11).Value = 17858660   # K33 (ALC): 16668361 -> 17858660
$ws.Cells.Item(33, 12).Value = 5420.2   # L33 (ALC): 5397.1816 -> 5420.2
$ws.Cells.Item(33, 13).Value = -17858431   # M33 (ALC): -16668132 -> -17858431
$ws.Cells.Item(33, 14).Value = -5878.2   # N33 (ALC): -5855.1816 -> -5878.2
$ws.Cells.Item(46, 8).Value = 2083.3333   # H46 (ALC): 2020 -> 2083.3333
$ws.Cells.Item(46, 10).Value = 2439.3333   # J46 (ALC): 2459 -> 2439.3333
$ws.Cells.Item(46, 12).Value = 7317.999899999999   # L46 (ALC): 7377 -> 7317.999899999999
$ws.Cells.Item(46, 14).Value = -7555.999899999999   # N46 (ALC): -7615 -> -7555.999899999999
$ws.Cells.Item(60, 8).Value = 2083.3333   # H60 (ALC): 2020 -> 2083.3333
$ws.Cells.Item(60, 10).Value = 2439.3333   # J60 (ALC): 2459 -> 2439.3333
$ws.Cells.Item(60, 12).Value = 7317.999899999999   # L60 (ALC): 7377 -> 7317.999899999999
$ws.Cells.Item(60, 14).Value = -8285.999899999999   # N60 (ALC): -8345 -> -8285.999899999999
$ws.Cells.Item(99, 8).Value = 1108.5   # H99 (ALC): 1047.091 -> 1108.5
$ws.Cells.Item(99, 9).Value = 987.44446   # I99 (ALC): 989.7778 -> 987.44446
$ws.Cells.Item(99, 10).Value = 1471.6666   # J99 (ALC): 1305 -> 1471.6666
$ws.Cells.Item(99, 11).Value = 2962.33338   # K99 (ALC): 2969.3334 -> 2962.33338
$ws.Cells.Item(99, 12).Value = 4414.9998   # L99 (ALC): 3915 -> 4414.9998
$ws.Cells.Item(99, 13).Value = -1464.33338   # M99 (ALC): -1471.3334 -> -1464.33338
$ws.Cells.Item(99, 14).Value = -7410.9998   # N99 (ALC): -6911 -> -7410.9998
$ws.Cells.Item(101, 8).Value = 1505.4445   # H101 (ALC): 1381.9 -> 1505.4445
$ws.Cells.Item(101, 9).Value = 290   # I101 (ALC): 283.33334 -> 290
$ws.Cells.Item(101, 11).Value = 870   # K101 (ALC): 850.0000200000001 -> 870
$ws.Cells.Item(101, 13).Value = 752   # M101 (ALC): 771.9999799999999 -> 752
$ws.Cells.Item(103, 8).Value = 494.2857   # H103 (ALC): 510.16666 -> 494.2857
$ws.Cells.Item(103, 10).Value = 411.33334   # J103 (ALC): 417.5 -> 411.33334
$ws.Cells.Item(103, 12).Value = 1234.00002   # L103 (ALC): 1252.5 -> 1234.00002
$ws.Cells.Item(103, 14).Value = -2406.00002   # N103 (ALC): -2424.5 -> -2406.00002
$ws.Cells.Item(116, 8).Value = 8799   # H116 (ALC): 9000 -> 8799
$ws.Cells.Item(116, 9).Value = 7998.3335   # I116 (ALC): 8000 -> 7998.3335
$ws.Cells.Item(116, 11).Value = 7998.3335   # K116 (ALC): 8000 -> 7998.3335
$ws.Cells.Item(116, 13).Value = -4556.3335   # M116 (ALC): -4558 -> -4556.3335
$ws.Cells.Item(123, 8).Value = 160000   # H123 (ALC): 250000 -> 160000
$ws.Cells.Item(123, 10).Value = 160000   # J123 (ALC): 250000 -> 160000
$ws.Cells.Item(123, 12).Value = 160000   # L123 (ALC): 250000 -> 160000
$ws.Cells.Item(123, 14).Value = -169800   # N123 (ALC): -259800 -> -169800
$ws.Cells.Item(124, 8).Value = 176666.67   # H124 (ALC): 220000 -> 176666.67
$ws.Cells.Item(124, 10).Value = 176666.67   # J124 (ALC): 220000 -> 176666.67
$ws.Cells.Item(124, 12).Value = 176666.67   # L124 (ALC): 220000 -> 176666.67
$ws.Cells.Item(124, 14).Value = -186486.67   # N124 (ALC): -229820 -> -186486.67
$ws.Cells.Item(125, 8).Value = 3500382.8   # H125 (ALC): 4136471.2 -> 3500382.8
$ws.Cells.Item(125, 10).Value = 1841   # J125 (ALC): 1787 -> 1841
$ws.Cells.Item(125, 12).Value = 16569   # L125 (ALC): 16083 -> 16569
$ws.Cells.Item(125, 14).Value = -21489   # N125 (ALC): -21003 -> -21489
$ws.Cells.Item(126, 8).Value = 172500   # H126 (ALC): 250000 -> 172500
$ws.Cells.Item(126, 10).Value = 172500   # J126 (ALC): 250000 -> 172500
$ws.Cells.Item(126, 12).Value = 172500   # L126 (ALC): 250000 -> 172500
$ws.Cells.Item(126, 14).Value = -182380   # N126 (ALC): -259880 -> -182380
$ws.Cells.Item(128, 8).Value = 250000   # H128 (ALC): 0 -> 250000
$ws.Cells.Item(128, 10).Value = 250000   # J128 (ALC): 0 -> 250000
$ws.Cells.Item(128, 12).Value = 250000   # L128 (ALC): 0 -> 250000
$ws.Cells.Item(128, 14).Value = -259960   # N128 (ALC): None -> -259960
$ws.Cells.Item(131, 8).Value = 2054.7144   # H131 (ALC): 2087.5715 -> 2054.7144
$ws.Cells.Item(131, 9).Value = 1064.6666   # I131 (ALC): 1103 -> 1064.6666
$ws.Cells.Item(131, 11).Value = 3193.9998   # K131 (ALC): 3309 -> 3193.9998
$ws.Cells.Item(131, 13).Value = 1846.0002   # M131 (ALC): 1731 -> 1846.0002
$ws.Cells.Item(135, 8).Value = 1924   # H135 (ALC): 2483.818 -> 1924
$ws.Cells.Item(135, 9).Value = 2150.7693   # I135 (ALC): 2483.818 -> 2150.7693
$ws.Cells.Item(135, 10).Value = 450   # J135 (ALC): 0 -> 450
$ws.Cells.Item(135, 11).Value = 19356.9237   # K135 (ALC): 22354.362 -> 19356.9237
$ws.Cells.Item(135, 12).Value = 4050   # L135 (ALC): 0 -> 4050
$ws.Cells.Item(135, 13).Value = -16821.9237   # M135 (ALC): -19819.362 -> -16821.9237
$ws.Cells.Item(135, 14).Value = -9120   # N135 (ALC): None -> -9120
$ws.Cells.Item(137, 8).Value = 5841.0435   # H137 (ALC): 5630.9585 -> 5841.0435
$ws.Cells.Item(137, 10).Value = 43200.2   # J137 (ALC): 31085.572 -> 43200.2
$ws.Cells.Item(137, 12).Value = 129600.6   # L137 (ALC): 93256.716 -> 129600.6
$ws.Cells.Item(137, 14).Value = -134700.6   # N137 (ALC): -98356.716 -> -134700.6
# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 4719636   # H32 (ARM): 4632400.5 -> 4719636
$ws.Cells.Item(32, 10).Value = 15665.667   # J32 (ARM): 13974.25 -> 15665.667
$ws.Cells.Item(32, 12).Value = 15665.667   # L32 (ARM): 13974.25 -> 15665.667
$ws.Cells.Item(32, 14).Value = -16239.667   # N32 (ARM): -14548.25 -> -16239.667
$ws.Cells.Item(60, 8).Value = 205013   # H60 (ARM): 188332.67 -> 205013
$ws.Cells.Item(60, 9).Value = 499995   # I60 (ARM): 269979.66 -> 499995
$ws.Cells.Item(60, 11).Value = 499995   # K60 (ARM): 269979.66 -> 499995
$ws.Cells.Item(60, 13).Value = -499262   # M60 (ARM): -269246.66 -> -499262
$ws.Cells.Item(63, 8).Value = 3399.6667   # H63 (ARM): 2585.4285 -> 3399.6667
$ws.Cells.Item(63, 9).Value = 3399.6667   # I63 (ARM): 2585.4285 -> 3399.6667
$ws.Cells.Item(63, 11).Value = 3399.6667   # K63 (ARM): 2585.4285 -> 3399.6667
$ws.Cells.Item(63, 13).Value = -2713.6667   # M63 (ARM): -1899.4285 -> -2713.6667
$ws.Cells.Item(66, 8).Value = 3399.6667   # H66 (ARM): 2585.4285 -> 3399.6667
$ws.Cells.Item(66, 9).Value = 3399.6667   # I66 (ARM): 2585.4285 -> 3399.6667
$ws.Cells.Item(66, 11).Value = 16998.3335   # K66 (ARM): 12927.1425 -> 16998.3335
$ws.Cells.Item(66, 13).Value = -13566.3335   # M66 (ARM): -9495.1425 -> -13566.3335
$ws.Cells.Item(110, 8).Value = 2757196.8   # H110 (ARM): 2843351.2 -> 2757196.8
$ws.Cells.Item(110, 9).Value = 3032836   # I110 (ARM): 3137408 -> 3032836
$ws.Cells.Item(110, 11).Value = 3032836   # K110 (ARM): 3137408 -> 3032836
$ws.Cells.Item(110, 13).Value = -3030791   # M110 (ARM): -3135363 -> -3030791
# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(94, 8).Value = 2304.5   # H94 (BSM): 2472.8333 -> 2304.5
$ws.Cells.Item(94, 9).Value = 1038.091   # I94 (BSM): 1113.45 -> 1038.091
$ws.Cells.Item(94, 11).Value = 1038.091   # K94 (BSM): 1113.45 -> 1038.091
$ws.Cells.Item(94, 13).Value = -587.0909999999999   # M94 (BSM): -662.45 -> -587.0909999999999
# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(6, 8).Value = 3631935.8   # H6 (CRP): 1667931.9 -> 3631935.8
$ws.Cells.Item(6, 9).Value = 4539420   # I6 (CRP): 3335033.2 -> 4539420
$ws.Cells.Item(6, 10).Value = 1999   # J6 (CRP): 830.3333 -> 1999
$ws.Cells.Item(6, 11).Value = 4539420   # K6 (CRP): 3335033.2 -> 4539420
$ws.Cells.Item(6, 12).Value = 1999   # L6 (CRP): 830.3333 -> 1999
$ws.Cells.Item(6, 13).Value = -4539307   # M6 (CRP): -3334920.2 -> -4539307
$ws.Cells.Item(6, 14).Value = -2225   # N6 (CRP): -1056.3333 -> -2225
$ws.Cells.Item(16, 8).Value = 8910.214   # H16 (CRP): 8018.3125 -> 8910.214
$ws.Cells.Item(16, 9).Value = 7505.75   # I16 (CRP): 6844 -> 7505.75
$ws.Cells.Item(16, 10).Value = 10782.833   # J16 (CRP): 9528.143 -> 10782.833
$ws.Cells.Item(16, 11).Value = 7505.75   # K16 (CRP): 6844 -> 7505.75
$ws.Cells.Item(16, 12).Value = 10782.833   # L16 (CRP): 9528.143 -> 10782.833
$ws.Cells.Item(16, 13).Value = -7218.75   # M16 (CRP): -6557 -> -7218.75
$ws.Cells.Item(16, 14).Value = -11356.833   # N16 (CRP): -10102.143 -> -11356.833
$ws.Cells.Item(19, 8).Value = 207.33333   # H19 (CRP): 2467.1667 -> 207.33333
$ws.Cells.Item(19, 9).Value = 61   # I19 (CRP): 2860.6 -> 61
$ws.Cells.Item(19, 11).Value = 61   # K19 (CRP): 2860.6 -> 61
$ws.Cells.Item(19, 13).Value = 109   # M19 (CRP): -2690.6 -> 109
$ws.Cells.Item(22, 8).Value = 941.8   # H22 (CRP): 1092.75 -> 941.8
$ws.Cells.Item(22, 9).Value = 525.25   # I22 (CRP): 600 -> 525.25
$ws.Cells.Item(22, 10).Value = 1219.5   # J22 (CRP): 1388.4 -> 1219.5
$ws.Cells.Item(22, 11).Value = 525.25   # K22 (CRP): 600 -> 525.25
$ws.Cells.Item(22, 12).Value = 1219.5   # L22 (CRP): 1388.4 -> 1219.5
$ws.Cells.Item(22, 13).Value = -175.25   # M22 (CRP): -250 -> -175.25
$ws.Cells.Item(22, 14).Value = -1919.5   # N22 (CRP): -2088.4 -> -1919.5
$ws.Cells.Item(24, 8).Value = 207.33333   # H24 (CRP): 2467.1667 -> 207.33333
$ws.Cells.Item(24, 9).Value = 61   # I24 (CRP): 2860.6 -> 61
$ws.Cells.Item(24, 11).Value = 61   # K24 (CRP): 2860.6 -> 61
$ws.Cells.Item(24, 13).Value = 109   # M24 (CRP): -2690.6 -> 109
$ws.Cells.Item(95, 8).Value = 18331.666   # H95 (CRP): 15000 -> 18331.666
$ws.Cells.Item(95, 10).Value = 18331.666   # J95 (CRP): 15000 -> 18331.666
$ws.Cells.Item(95, 12).Value = 18331.666   # L95 (CRP): 15000 -> 18331.666
$ws.Cells.Item(95, 14).Value = -23823.666   # N95 (CRP): -20492 -> -23823.666
$ws.Cells.Item(110, 8).Value = 105000   # H110 (CRP): 0 -> 105000
$ws.Cells.Item(110, 10).Value = 105000   # J110 (CRP): 0 -> 105000
$ws.Cells.Item(110, 12).Value = 105000   # L110 (CRP): 0 -> 105000
$ws.Cells.Item(110, 14).Value = -113180   # N110 (CRP): None -> -113180
$ws.Cells.Item(113, 8).Value = 8910.214   # H113 (CRP): 8018.3125 -> 8910.214
$ws.Cells.Item(113, 9).Value = 7505.75   # I113 (CRP): 6844 -> 7505.75
$ws.Cells.Item(113, 10).Value = 10782.833   # J113 (CRP): 9528.143 -> 10782.833
$ws.Cells.Item(113, 11).Value = 7505.75   # K113 (CRP): 6844 -> 7505.75
$ws.Cells.Item(113, 12).Value = 10782.833   # L113 (CRP): 9528.143 -> 10782.833
$ws.Cells.Item(113, 13).Value = -5335.75   # M113 (CRP): -4674 -> -5335.75
$ws.Cells.Item(113, 14).Value = -15122.833   # N113 (CRP): -13868.143 -> -15122.833
# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(80, 8).Value = 10775   # H80 (CUL): 10791.5 -> 10775
$ws.Cells.Item(80, 9).Value = 3966.6667   # I80 (CUL): 3999.6667 -> 3966.6667
$ws.Cells.Item(80, 11).Value = 11900.0001   # K80 (CUL): 11999.0001 -> 11900.0001
$ws.Cells.Item(80, 13).Value = -10964.0001   # M80 (CUL): -11063.0001 -> -10964.0001
$ws.Cells.Item(83, 8).Value = 10775   # H83 (CUL): 10791.5 -> 10775
$ws.Cells.Item(83, 9).Value = 3966.6667   # I83 (CUL): 3999.6667 -> 3966.6667
$ws.Cells.Item(83, 11).Value = 35700.0003   # K83 (CUL): 35997.0003 -> 35700.0003
$ws.Cells.Item(83, 13).Value = -31020.0003   # M83 (CUL): -31317.0003 -> -31020.0003
$ws.Cells.Item(117, 8).Value = 598.7778   # H117 (CUL): 685.7143 -> 598.7778
$ws.Cells.Item(117, 9).Value = 372.25   # I117 (CUL): 450 -> 372.25
$ws.Cells.Item(117, 11).Value = 1116.75   # K117 (CUL): 1350 -> 1116.75
$ws.Cells.Item(117, 13).Value = 2325.25   # M117 (CUL): 2092 -> 2325.25
$ws.Cells.Item(131, 8).Value = 1445.61   # H131 (CUL): 1441.85 -> 1445.61
$ws.Cells.Item(131, 9).Value = 866.6667   # I131 (CUL): 849.5 -> 866.6667
$ws.Cells.Item(131, 10).Value = 1463.5155   # J131 (CUL): 1466.5312 -> 1463.5155
$ws.Cells.Item(131, 11).Value = 2600.0001   # K131 (CUL): 2548.5 -> 2600.0001
$ws.Cells.Item(131, 12).Value = 4390.5465   # L131 (CUL): 4399.5936 -> 4390.5465
$ws.Cells.Item(131, 13).Value = 2439.9999   # M131 (CUL): 2491.5 -> 2439.9999
$ws.Cells.Item(131, 14).Value = -14470.5465   # N131 (CUL): -14479.5936 -> -14470.5465
$ws.Cells.Item(139, 8).Value = 10563.556   # H139 (CUL): 8081.2915 -> 10563.556
$ws.Cells.Item(139, 9).Value = 13553.75   # I139 (CUL): 8971.157999999999 -> 13553.75
$ws.Cells.Item(139, 10).Value = 4583.1665   # J139 (CUL): 4699.8 -> 4583.1665
$ws.Cells.Item(139, 11).Value = 40661.25   # K139 (CUL): 26913.474 -> 40661.25
$ws.Cells.Item(139, 12).Value = 13749.4995   # L139 (CUL): 14099.4 -> 13749.4995
$ws.Cells.Item(139, 13).Value = -35521.25   # M139 (CUL): -21773.474 -> -35521.25
$ws.Cells.Item(139, 14).Value = -24029.4995   # N139 (CUL): -24379.4 -> -24029.4995
# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 8).Value = 149.75   # H2 (GSM): 147.72 -> 149.75
$ws.Cells.Item(2, 9).Value = 122.933334   # I2 (GSM): 121.4375 -> 122.933334
$ws.Cells.Item(2, 11).Value = 122.933334   # K2 (GSM): 121.4375 -> 122.933334
$ws.Cells.Item(2, 13).Value = -9.933334000000002   # M2 (GSM): -8.4375 -> -9.933334000000002
$ws.Cells.Item(5, 8).Value = 5999.5   # H5 (GSM): 5000 -> 5999.5
$ws.Cells.Item(5, 10).Value = 6999   # J5 (GSM): 0 -> 6999
$ws.Cells.Item(5, 12).Value = 6999   # L5 (GSM): 0 -> 6999
$ws.Cells.Item(5, 14).Value = -7223   # N5 (GSM): None -> -7223
$ws.Cells.Item(7, 8).Value = 500500   # H7 (GSM): 525000 -> 500500
$ws.Cells.Item(7, 9).Value = 1000   # I7 (GSM): 50000 -> 1000
$ws.Cells.Item(7, 11).Value = 1000   # K7 (GSM): 50000 -> 1000
$ws.Cells.Item(7, 13).Value = -888   # M7 (GSM): -49888 -> -888
$ws.Cells.Item(8, 8).Value = 500500   # H8 (GSM): 525000 -> 500500
$ws.Cells.Item(8, 9).Value = 1000   # I8 (GSM): 50000 -> 1000
$ws.Cells.Item(8, 11).Value = 1000   # K8 (GSM): 50000 -> 1000
$ws.Cells.Item(8, 13).Value = -861   # M8 (GSM): -49861 -> -861
$ws.Cells.Item(11, 8).Value = 576237.2   # H11 (GSM): 560316.75 -> 576237.2
$ws.Cells.Item(11, 9).Value = 692444.8   # I11 (GSM): 669466.7 -> 692444.8
$ws.Cells.Item(11, 11).Value = 692444.8   # K11 (GSM): 669466.7 -> 692444.8
$ws.Cells.Item(11, 13).Value = -692305.8   # M11 (GSM): -669327.7 -> -692305.8
$ws.Cells.Item(12, 8).Value = 5001.3335   # H12 (GSM): 4976 -> 5001.3335
$ws.Cells.Item(12, 10).Value = 5001.3335   # J12 (GSM): 4976 -> 5001.3335
$ws.Cells.Item(12, 12).Value = 5001.3335   # L12 (GSM): 4976 -> 5001.3335
$ws.Cells.Item(12, 14).Value = -5281.3335   # N12 (GSM): -5256 -> -5281.3335
$ws.Cells.Item(58, 8).Value = 29000   # H58 (GSM): 27497.5 -> 29000
$ws.Cells.Item(58, 9).Value = 0   # I58 (GSM): 25995 -> 0
$ws.Cells.Item(58, 11).Value = 0   # K58 (GSM): 25995 -> 0
$ws.Cells.Item(58, 13).Value = $null   # M58 (GSM): -25718 -> None
$ws.Cells.Item(70, 8).Value = 6520   # H70 (GSM): 6524.35 -> 6520
$ws.Cells.Item(70, 9).Value = 5082.3335   # I70 (GSM): 5089.5835 -> 5082.3335
$ws.Cells.Item(70, 11).Value = 5082.3335   # K70 (GSM): 5089.5835 -> 5082.3335
$ws.Cells.Item(70, 13).Value = -4812.3335   # M70 (GSM): -4819.5835 -> -4812.3335
$ws.Cells.Item(73, 8).Value = 6520   # H73 (GSM): 6524.35 -> 6520
$ws.Cells.Item(73, 9).Value = 5082.3335   # I73 (GSM): 5089.5835 -> 5082.3335
$ws.Cells.Item(73, 11).Value = 5082.3335   # K73 (GSM): 5089.5835 -> 5082.3335
$ws.Cells.Item(73, 13).Value = -4146.3335   # M73 (GSM): -4153.5835 -> -4146.3335
$ws.Cells.Item(80, 8).Value = 13007.625   # H80 (GSM): 15787.827 -> 13007.625
$ws.Cells.Item(80, 9).Value = 10408.125   # I80 (GSM): 15066.571 -> 10408.125
$ws.Cells.Item(80, 10).Value = 15607.125   # J80 (GSM): 16461 -> 15607.125
$ws.Cells.Item(80, 11).Value = 10408.125   # K80 (GSM): 15066.571 -> 10408.125
$ws.Cells.Item(80, 12).Value = 15607.125   # L80 (GSM): 16461 -> 15607.125
$ws.Cells.Item(80, 13).Value = -9410.125   # M80 (GSM): -14068.571 -> -9410.125
$ws.Cells.Item(80, 14).Value = -17603.125   # N80 (GSM): -18457 -> -17603.125
$ws.Cells.Item(83, 8).Value = 13007.625   # H83 (GSM): 15787.827 -> 13007.625
$ws.Cells.Item(83, 9).Value = 10408.125   # I83 (GSM): 15066.571 -> 10408.125
$ws.Cells.Item(83, 10).Value = 15607.125   # J83 (GSM): 16461 -> 15607.125
$ws.Cells.Item(83, 11).Value = 52040.625   # K83 (GSM): 75332.855 -> 52040.625
$ws.Cells.Item(83, 12).Value = 78035.625   # L83 (GSM): 82305 -> 78035.625
$ws.Cells.Item(83, 13).Value = -47048.625   # M83 (GSM): -70340.855 -> -47048.625
$ws.Cells.Item(83, 14).Value = -88019.625   # N83 (GSM): -92289 -> -88019.625
$ws.Cells.Item(122, 8).Value = 1702097.4   # H122 (GSM): 1261409.6 -> 1702097.4
$ws.Cells.Item(122, 9).Value = 1791565.2   # I122 (GSM): 1547573.1 -> 1791565.2
$ws.Cells.Item(122, 10).Value = 2208   # J122 (GSM): 2290 -> 2208
$ws.Cells.Item(122, 11).Value = 5374695.6   # K122 (GSM): 4642719.300000001 -> 5374695.6
$ws.Cells.Item(122, 12).Value = 6624   # L122 (GSM): 6870 -> 6624
$ws.Cells.Item(122, 13).Value = -5372245.6   # M122 (GSM): -4640269.300000001 -> -5372245.6
$ws.Cells.Item(122, 14).Value = -11524   # N122 (GSM): -11770 -> -11524
# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(2, 8).Value = 0   # H2 (LTW): 1 -> 0
$ws.Cells.Item(2, 9).Value = 0   # I2 (LTW): 1 -> 0
$ws.Cells.Item(2, 11).Value = 0   # K2 (LTW): 1 -> 0
$ws.Cells.Item(2, 13).Value = $null   # M2 (LTW): 111 -> None
$ws.Cells.Item(11, 8).Value = 9998   # H11 (LTW): 15426.625 -> 9998
$ws.Cells.Item(11, 10).Value = 9998   # J11 (LTW): 15426.625 -> 9998
$ws.Cells.Item(11, 12).Value = 9998   # L11 (LTW): 15426.625 -> 9998
$ws.Cells.Item(11, 14).Value = -10278   # N11 (LTW): -15706.625 -> -10278
$ws.Cells.Item(38, 8).Value = 105199   # H38 (LTW): 100165.664 -> 105199
$ws.Cells.Item(38, 10).Value = 105199   # J38 (LTW): 100165.664 -> 105199
$ws.Cells.Item(38, 12).Value = 105199   # L38 (LTW): 100165.664 -> 105199
$ws.Cells.Item(38, 14).Value = -106019   # N38 (LTW): -100985.664 -> -106019
$ws.Cells.Item(40, 8).Value = 1594246.1   # H40 (LTW): 1638494.6 -> 1594246.1
$ws.Cells.Item(40, 9).Value = 3302.8333   # I40 (LTW): 3391.2942 -> 3302.8333
$ws.Cells.Item(40, 10).Value = 3101455.5   # J40 (LTW): 3101482 -> 3101455.5
$ws.Cells.Item(40, 11).Value = 3302.8333   # K40 (LTW): 3391.2942 -> 3302.8333
$ws.Cells.Item(40, 12).Value = 3101455.5   # L40 (LTW): 3101482 -> 3101455.5
$ws.Cells.Item(40, 13).Value = -3166.8333   # M40 (LTW): -3255.2942 -> -3166.8333
$ws.Cells.Item(40, 14).Value = -3101727.5   # N40 (LTW): -3101754 -> -3101727.5
$ws.Cells.Item(46, 8).Value = 1918.6957   # H46 (LTW): 1862.2916 -> 1918.6957
$ws.Cells.Item(46, 9).Value = 547.25   # I46 (LTW): 546.25 -> 547.25
$ws.Cells.Item(46, 10).Value = 2207.4211   # J46 (LTW): 2125.5 -> 2207.4211
$ws.Cells.Item(46, 11).Value = 547.25   # K46 (LTW): 546.25 -> 547.25
$ws.Cells.Item(46, 12).Value = 2207.4211   # L46 (LTW): 2125.5 -> 2207.4211
$ws.Cells.Item(46, 13).Value = -359.25   # M46 (LTW): -358.25 -> -359.25
$ws.Cells.Item(46, 14).Value = -2583.4211   # N46 (LTW): -2501.5 -> -2583.4211
$ws.Cells.Item(55, 8).Value = 1698.8334   # H55 (LTW): 1703.4166 -> 1698.8334
$ws.Cells.Item(55, 9).Value = 1279.5   # I55 (LTW): 1366.3846 -> 1279.5
$ws.Cells.Item(55, 10).Value = 2285.9   # J55 (LTW): 2101.7273 -> 2285.9
$ws.Cells.Item(55, 11).Value = 1279.5   # K55 (LTW): 1366.3846 -> 1279.5
$ws.Cells.Item(55, 12).Value = 2285.9   # L55 (LTW): 2101.7273 -> 2285.9
$ws.Cells.Item(55, 13).Value = -1106.5   # M55 (LTW): -1193.3846 -> -1106.5
$ws.Cells.Item(55, 14).Value = -2631.9   # N55 (LTW): -2447.7273 -> -2631.9
$ws.Cells.Item(58, 8).Value = 3697.6667   # H58 (LTW): 4807.4 -> 3697.6667
$ws.Cells.Item(58, 9).Value = 1193   # I58 (LTW): 3796.75 -> 1193
$ws.Cells.Item(58, 10).Value = 4950   # J58 (LTW): 8850 -> 4950
$ws.Cells.Item(58, 11).Value = 1193   # K58 (LTW): 3796.75 -> 1193
$ws.Cells.Item(58, 12).Value = 4950   # L58 (LTW): 8850 -> 4950
$ws.Cells.Item(58, 13).Value = -933   # M58 (LTW): -3536.75 -> -933
$ws.Cells.Item(58, 14).Value = -5470   # N58 (LTW): -9370 -> -5470
$ws.Cells.Item(136, 8).Value = 1116342.8   # H136 (LTW): 1076480.5 -> 1116342.8
$ws.Cells.Item(136, 9).Value = 16722.77   # I136 (LTW): 15542.571 -> 16722.77
$ws.Cells.Item(136, 11).Value = 50168.31   # K136 (LTW): 46627.713 -> 50168.31
$ws.Cells.Item(136, 13).Value = -47618.31   # M136 (LTW): -44077.713 -> -47618.31
# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(2, 8).Value = 591305.4399999999   # H2 (WVR): 772696.4399999999 -> 591305.4399999999
$ws.Cells.Item(2, 9).Value = 772961.7   # I2 (WVR): 836837.8 -> 772961.7
$ws.Cells.Item(2, 10).Value = 922.5   # J2 (WVR): 3000 -> 922.5
$ws.Cells.Item(2, 11).Value = 772961.7   # K2 (WVR): 836837.8 -> 772961.7
$ws.Cells.Item(2, 12).Value = 922.5   # L2 (WVR): 3000 -> 922.5
$ws.Cells.Item(2, 13).Value = -772849.7   # M2 (WVR): -836725.8 -> -772849.7
$ws.Cells.Item(2, 14).Value = -1146.5   # N2 (WVR): -3224 -> -1146.5
$ws.Cells.Item(4, 8).Value = 1053243.1   # H4 (WVR): 737.2353000000001 -> 1053243.1
$ws.Cells.Item(4, 9).Value = 388.18182   # I4 (WVR): 430 -> 388.18182
$ws.Cells.Item(4, 10).Value = 2500918.8   # J4 (WVR): 1300.5 -> 2500918.8
$ws.Cells.Item(4, 11).Value = 388.18182   # K4 (WVR): 430 -> 388.18182
$ws.Cells.Item(4, 12).Value = 2500918.8   # L4 (WVR): 1300.5 -> 2500918.8
$ws.Cells.Item(4, 13).Value = -275.18182   # M4 (WVR): -317 -> -275.18182
$ws.Cells.Item(4, 14).Value = -2501144.8   # N4 (WVR): -1526.5 -> -2501144.8
$ws.Cells.Item(5, 8).Value = 4498.3335   # H5 (WVR): 7495 -> 4498.3335
$ws.Cells.Item(5, 9).Value = 4000   # I5 (WVR): 0 -> 4000
$ws.Cells.Item(5, 10).Value = 4747.5   # J5 (WVR): 7495 -> 4747.5
$ws.Cells.Item(5, 11).Value = 4000   # K5 (WVR): 0 -> 4000
$ws.Cells.Item(5, 12).Value = 4747.5   # L5 (WVR): 7495 -> 4747.5
$ws.Cells.Item(5, 13).Value = -3888   # M5 (WVR): None -> -3888
$ws.Cells.Item(5, 14).Value = -4971.5   # N5 (WVR): -7719 -> -4971.5
$ws.Cells.Item(9, 8).Value = 87989.8   # H9 (WVR): 219184.5 -> 87989.8
$ws.Cells.Item(9, 9).Value = 87989.8   # I9 (WVR): 219184.5 -> 87989.8
$ws.Cells.Item(9, 11).Value = 87989.8   # K9 (WVR): 219184.5 -> 87989.8
$ws.Cells.Item(9, 13).Value = -87849.8   # M9 (WVR): -219044.5 -> -87849.8
$ws.Cells.Item(14, 8).Value = 35334.668   # H14 (WVR): 26126.75 -> 35334.668
$ws.Cells.Item(14, 9).Value = 3002   # I14 (WVR): 1003.5 -> 3002
$ws.Cells.Item(14, 10).Value = 100000   # J14 (WVR): 51250 -> 100000
$ws.Cells.Item(14, 11).Value = 3002   # K14 (WVR): 1003.5 -> 3002
$ws.Cells.Item(14, 12).Value = 100000   # L14 (WVR): 51250 -> 100000
$ws.Cells.Item(14, 13).Value = -2834   # M14 (WVR): -835.5 -> -2834
$ws.Cells.Item(14, 14).Value = -100336   # N14 (WVR): -51586 -> -100336
$ws.Cells.Item(81, 8).Value = 11251   # H81 (WVR): 20501.5 -> 11251
$ws.Cells.Item(81, 9).Value = 14334   # I81 (WVR): 39001 -> 14334
$ws.Cells.Item(81, 11).Value = 28668   # K81 (WVR): 78002 -> 28668
$ws.Cells.Item(81, 13).Value = -27607   # M81 (WVR): -76941 -> -27607
$ws.Cells.Item(84, 8).Value = 11251   # H84 (WVR): 20501.5 -> 11251
$ws.Cells.Item(84, 9).Value = 14334   # I84 (WVR): 39001 -> 14334
$ws.Cells.Item(84, 11).Value = 143340   # K84 (WVR): 390010 -> 143340
$ws.Cells.Item(84, 13).Value = -138036   # M84 (WVR): -384706 -> -138036
$ws.Cells.Item(113, 8).Value = 3291.423   # H113 (WVR): 3403.44 -> 3291.423
$ws.Cells.Item(113, 9).Value = 3360.0952   # I113 (WVR): 3503.55 -> 3360.0952
$ws.Cells.Item(113, 11).Value = 10080.2856   # K113 (WVR): 10510.65 -> 10080.2856
$ws.Cells.Item(113, 13).Value = -7910.285600000001   # M113 (WVR): -8340.650000000001 -> -7910.285600000001
$ws.Cells.Item(122, 8).Value = 414490.47   # H122 (WVR): 422361.78 -> 414490.47
$ws.Cells.Item(122, 9).Value = 571742.8   # I122 (WVR): 586859.5 -> 571742.8
$ws.Cells.Item(122, 11).Value = 1715228.4   # K122 (WVR): 1760578.5 -> 1715228.4
$ws.Cells.Item(122, 13).Value = -1712778.4   # M122 (WVR): -1758128.5 -> -1712778.4
$ws.Cells.Item(132, 8).Value = 394419.78   # H132 (WVR): 409006.12 -> 394419.78
$ws.Cells.Item(132, 9).Value = 3151.8572   # I132 (WVR): 3280 -> 3151.8572
$ws.Cells.Item(132, 11).Value = 9455.571599999999   # K132 (WVR): 9840 -> 9455.571599999999
$ws.Cells.Item(132, 13).Value = -6925.571599999999   # M132 (WVR): -7310 -> -6925.571599999999
